$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing
# header cell H1 (bold font, centered alignment, thin border) by copying
# its format (reuses the same cell style rather than creating a new one).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new I and J columns (rows 2-24)
$values = @(
    @(6, 6),
    @(9, 9),
    @(6, 8),
    @(5, 6),
    @(6, 6),
    @(8, 8),
    @(11, 11),
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(4, 5),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(4, 4)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
